# Update pl_mw line-loss results for the 380 kV case (rows 2-25, columns B,D,E,F,G,I,J,L,O)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","D","E","F","G","I","J","L","O")

$data = @(
    @(1.67961083434102, 0.1937505428121185, 0.2083201996281332, 1.213919806098907, 0.002437477770410056, 0.8346540536711338, 0.2845129632374039, 0.8194605269708006, 3.003273676669608),
    @(1.540912147450399, 0.1935376105644195, 0.203688917155592, 1.225109026517927, 0.002440928170151561, 0.8488667628801529, 0.2732944691124004, 0.738676726841021, 3.009311621701244),
    @(1.455552984286669, 0.1934513796467527, 0.2008915190453138, 1.232960815413492, 0.002443160529698022, 0.8583162854298259, 0.2664586793000154, 0.6889303194444381, 3.015558088584612),
    @(1.420721070384957, 0.1934275086641293, 0.1997633776144738, 1.236406793098801, 0.00244409893871912, 0.8623486505704037, 0.2636865911503747, 0.6686232541694892, 3.018740414737067),
    @(1.414934466330408, 0.1934242277652416, 0.1995767703638549, 1.236993856568866, 0.002444256497156046, 0.8630291858069477, 0.263227117546009, 0.6652492096013702, 3.019307252853451),
    @(1.45508341757369, 0.1934510119850756, 0.2008762564107442, 1.233006292458001, 0.002443173069159115, 0.8583699322045177, 0.266421238582609, 0.6886565910809281, 3.015598430205273),
    @(1.631829884199988, 0.1936679309695251, 0.2067138544160585, 1.217573817647875, 0.002438643902119786, 0.8394044673867818, 0.2806341669305112, 0.7916371520637711, 3.004827561442255),
    @(1.976777137161434, 0.1944432645404959, 0.2185202691975761, 1.195122669839229, 0.002430661141374427, 0.8079563668307692, 0.3089068205727727, 0.9923806435823224, 3.003931274173226),
    @(2.229116518662124, 0.1952221294864174, 0.2274037122669554, 1.183423649672335, 0.002425338563347822, 0.7883633986378413, 0.3299059395533561, 1.139080109274914, 3.015717018884914),
    @(2.343658760216329, 0.1956210349227732, 0.2314885383748759, 1.17915008735639, 0.002423033748723961, 0.7802151206256553, 0.3395048534173384, 1.205636303294739, 3.02380574558947),
    @(2.386995383231579, 0.1957784295239904, 0.2330414617195089, 1.177683134857233, 0.002422177629837652, 0.7772397640082644, 0.3431460408939415, 1.230812625407452, 3.027262830617246),
    @(2.377663794759087, 0.1957442510957748, 0.232706744211022, 1.177992328375424, 0.002422361270588739, 0.777875654007012, 0.3423615722757205, 1.225391680528276, 3.026500729792758),
    @(2.347224865910334, 0.1956338572806331, 0.2316161774092151, 1.179026363541226, 0.00242296298180257, 0.7799681257592539, 0.3398042922141826, 1.207708125103295, 3.024082253794148),
    @(2.328575120457344, 0.1955670610289957, 0.2309489608524444, 1.179679468234752, 0.002423333715145182, 0.7812641868434298, 0.3382386921406209, 1.196872869582307, 3.022652241159079),
    @(2.221625717782388, 0.1951969510314413, 0.2271376218810488, 1.183724084175765, 0.002425491524770818, 0.7889113169545681, 0.3292795279044753, 1.13472679670798, 3.015243442618328),
    @(2.155950469634433, 0.1949812690135602, 0.2248105420959874, 1.186474265629478, 0.002426845039458891, 0.793798615991502, 0.3237949655203209, 1.096555533230799, 3.011398253302815),
    @(2.118152579084892, 0.1948614164453062, 0.2234761925002786, 1.188154740192928, 0.002427634511046702, 0.7966816185776295, 0.3206447710953881, 1.074583751127534, 3.009443199734164),
    @(2.105350946033127, 0.194821560423641, 0.2230251191642552, 1.188740643820786, 0.002427903698625612, 0.7976701046100061, 0.3195789332122274, 1.067141662755489, 3.008825265987525),
    @(2.162944132714586, 0.1950037943838012, 0.2250578382019697, 1.186171290613025, 0.002426699821209435, 0.7932709054613021, 0.3243783558631037, 1.100620666372095, 3.011781009792799),
    @(2.356166568005733, 0.1956661112101443, 0.2319363396309697, 1.178718529924566, 0.002422785793004145, 0.7793505226372091, 0.3405552596693013, 1.212902959951521, 3.024781909238527),
    @(2.482225823628767, 0.1961358804965485, 0.2364672360151729, 1.174730196121899, 0.002420324844953747, 0.7708953173846353, 0.3511642782089126, 1.286127520767423, 3.035576497769625),
    @(2.414966775524704, 0.1958818034477332, 0.234045836663725, 1.176777901390608, 0.002421629441760948, 0.7753491360320126, 0.3454988294725894, 1.247061162733132, 3.029604360443841),
    @(2.159782423029867, 0.1949935977449613, 0.2249460246544146, 1.186307956160988, 0.002426765439050719, 0.7935092553034089, 0.3241145958606921, 1.098782902476671, 3.011607169540582),
    @(1.883645695261407, 0.1941964786673651, 0.2152888913321078, 1.200356516219742, 0.002432725040187552, 0.8158479829234402, 0.3012171848083796, 0.9382085128125937, 3.001996652292263),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $values[$j]
    }
}

